# Update "想去人数" (want-to-go count) figures in the F column of the
# "展览" (index 1) and "全部类型" (index 4) sheets to the refreshed
# scrape values from the gh-pages data regeneration.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value  = 8710
$wsExhibit.Range("F6").Value  = 502
$wsExhibit.Range("F7").Value  = 182
$wsExhibit.Range("F8").Value  = 19
$wsExhibit.Range("F9").Value  = 481
$wsExhibit.Range("F11").Value = 95
$wsExhibit.Range("F13").Value = 6316
$wsExhibit.Range("F14").Value = 213
$wsExhibit.Range("F15").Value = 336
$wsExhibit.Range("F16").Value = 2466
$wsExhibit.Range("F17").Value = 133
$wsExhibit.Range("F18").Value = 228
$wsExhibit.Range("F19").Value = 250
$wsExhibit.Range("F20").Value = 482

# --- Sheet 4: 全部类型 ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 8710
$wsAll.Range("F8").Value  = 502
$wsAll.Range("F9").Value  = 182
$wsAll.Range("F10").Value = 19
$wsAll.Range("F11").Value = 481
$wsAll.Range("F13").Value = 95
$wsAll.Range("F16").Value = 6316
$wsAll.Range("F18").Value = 213
$wsAll.Range("F19").Value = 336
$wsAll.Range("F20").Value = 2466
$wsAll.Range("F21").Value = 133
$wsAll.Range("F22").Value = 228
$wsAll.Range("F23").Value = 250
$wsAll.Range("F24").Value = 482
